$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Quotité proposée" -> "Quotité"
$ws.Range("B1").Value = "Quotité"

# Full rebuild of A2:B72 (department list changed: Corrèze (19) and
# Eure-et-Loire (28) added as their own rows, Limousin/Loiret labels
# adjusted, all quotité values refreshed, 2 new rows appended before Total)
$names = @(
    'Ain (01)',
    'Aisne (02)',
    'Alpes de Haute-Provence (04)',
    'Alpes Maritimes (06)',
    'Alsace (67+68)',
    'Aquitaine (24+33+40+47+64)',
    'Ariège (09)',
    'Aveyron (12)',
    'Bouches du Rhône (13)',
    'Bourgogne (71+58+21)',
    'Calvados (14)',
    'Cantal (15)',
    'Champagne-Ardennes (51+08+52+10)',
    'Charente (16)',
    'Charente Maritime (17)',
    'Cher (18)',
    'Corrèze (19)',
    'Côtes d''Armor (22)',
    'Deux-Sèvres (79)',
    'Eure-et-Loire (28)',
    'Essonne (91)',
    'Finistère (29)',
    'Franche-Comté (25+39+70+90)',
    'Gard (30+48)',
    'Gers (32)',
    'Grenoble (07+26+38+73+74)',
    'Guyane (973)',
    'Haute-Garonne (31+65)',
    'Haute-Loire (43)',
    'Haute-Normandie (27+76)',
    'Hautes-Alpes (05)',
    'Hauts de Seine (92)',
    'Hérault (34)',
    'Ile et Vilaine (35)',
    'Indre (36)',
    'Indre et Loire (37)',
    'Limousin (23+87)',
    'Loir et Cher (41)',
    'Loire (42)',
    'Loire Atlantique (44)',
    'Loiret (45)',
    'Lorraine (54+55+57+88)',
    'Lot (46)',
    'Maine et Loire (49)',
    'Manche (50)',
    'Mayenne (53)',
    'Mayotte (976)',
    'Morbihan (56)',
    'Nord (59)',
    'Oise (60)',
    'Orne (61)',
    'Paris (75)',
    'Pas-de-Calais (62)',
    'Puy de Dôme - Allier (03+63)',
    'Pyrénées Orientales - Aude (66+11)',
    'Réunion (974)',
    'Rhône (69)',
    'Sarthe (72)',
    'Seine-et-Marne (77)',
    'Seine-Saint-Denis (93)',
    'Somme (80)',
    'Tarn (81)',
    'Val d''Oise (95)',
    'Val de Marne (94)',
    'Var (83)',
    'Vaucluse (84)',
    'Vendée (85)',
    'Vienne (86)',
    'Yonne (89)',
    'Yvelines (78)',
    'Total'
)

$values = @(
    0.567,
    0.588,
    0.567,
    0.466,
    0.948,
    0.899,
    0.49,
    0.575,
    1.087,
    0.729,
    1.156,
    0.466,
    0.765,
    0.616,
    0.571,
    0.462,
    0.478,
    0.73,
    0.575,
    0.478,
    0.953,
    1.018,
    0.932,
    0.749,
    0.494,
    1.42,
    0.806,
    1.119,
    0.506,
    0.993,
    0.62,
    0.965,
    1.27,
    0.754,
    0.445,
    0.518,
    0.676,
    0.449,
    0.697,
    0.884,
    0.757,
    1.046,
    0.506,
    0.669,
    0.896,
    0.547,
    0.604,
    0.811,
    1.075,
    0.506,
    0.575,
    1.9,
    0.783,
    1.481,
    0.741,
    0.555,
    1.392,
    0.518,
    0.787,
    2.221,
    0.555,
    0.571,
    0.774,
    1.051,
    0.685,
    0.547,
    0.754,
    0.583,
    0.559,
    0.575,
    54.505
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Match formatting of the two brand-new rows (71-72) to the rest of the table
$ws.Range("A71:B72").Style = $ws.Range("A70:B70").Style

# Restore dimension/selection state seen after editing (scrolled to bottom,
# last cell B72 selected)
$ws.Range("B72").Select() | Out-Null
